$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new meeting-minutes row (row 14), matching the formatting of the
# row above it (row 13) by copying formats first, then filling in values.
$ws.Range("A13:C13").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A14").Value = "November 2 2023"
$ws.Range("B14").Value = "Madison"
$ws.Range("C14").Value = "7:00PM"

# Update the view: active cell moves to A9 (scrolled back into view, so the
# previous topLeftCell override is no longer needed).
$ws.Range("A9").Select() | Out-Null
